$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "hQbLc655"
$ws.Range("B2").Value = 231004193
$ws.Range("C2").Value = "kkxbdqq58"
$ws.Range("D2").Value = "Hp%93Ny#"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "lBKZgNJu"
$ws.Range("G2").Value = "UTWa"
$ws.Range("H2").Value = "Candidate"

# Delete row 3 entirely (was removed from the sheet)
$ws.Rows("3:3").Delete()

$ws.Range("A1:H2").Select() | Out-Null
